# Actualizacion de la lista de bugs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: fill in bug #17 ---
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "Equipos para generar fixture"
$ws.Range("C19").Value = "Deberían ser como mínimo 2"
$ws.Range("D19").Value = "Flor"
$ws.Range("E19").Value = "edicion-configurar.aspx"
$ws.Range("F19").Value = "PENDIENTE"

# Copy the "vertical centered" text format (used by column E elsewhere, e.g. E4)
# onto E19, and the highlighted "PENDIENTE" status format (from F18) onto F19.
$ws.Range("E4").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("F18").Copy()
$ws.Range("F19").PasteSpecial(-4122)

# --- Row 20: fill in bug #18 ---
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "modificacion de configuracion de edicion"
$ws.Range("C20").Value = "cuando se modific la congiuracion de la edicion deberia generar el fixture con los nuevos equipos"
$ws.Range("D20").Value = "Flor"
$ws.Range("E20").Value = "edicion-configurar.aspx"
$ws.Range("F20").Value = "PENDIENTE"

# Copy the wrapped-text format (used by C17/C18) onto B20/C20, the vertical
# centered format onto E20, and the "PENDIENTE" status format onto F20.
$ws.Range("C18").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("F18").Copy()
$ws.Range("F20").PasteSpecial(-4122)

$ws.Rows.Item(20).RowHeight = 45

# --- Update the view: scrolled down with C17 selected ---
$ws.Range("C17").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1

$excel.CutCopyMode = $false
